$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.826.66"
$ws.Range("E2").Value = "  +4.37%  "
$ws.Range("D3").Value = "'2.275.09"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'302.92"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'101.44"
$ws.Range("E6").Value = "  +8.33%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").Value = "'35.85"
$ws.Range("E10").Value = "  +5.47%  "
$ws.Range("D11").Value = "'0.0785"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'2.623.63"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "'2.276.30"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "'13.71"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'46.827.96"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'0.801"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'0.0₃0935"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "'65.50"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'249.64"
$ws.Range("E23").Value = "  +4.73%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "'42.70"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'9.75"
$ws.Range("E29").Value = "  +2.45%  "
$ws.Range("D30").Value = "'19.93"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").Value = "'2.80"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("D32").Value = "'5.51"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "'147.43"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("E34").Value = "  +13.90%  "
$ws.Range("D35").Value = "'0.0779"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'0.115"
$ws.Range("E36").Value = "  +10.73%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "'16.15"
$ws.Range("E38").Value = "  +19.21%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'3.93"
$ws.Range("E40").Value = "  +3.97%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0299"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "'3.26"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "'1.809.80"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "'90.14"
$ws.Range("E46").Value = "  +19.83%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.191"
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "'72.95"
$ws.Range("E48").Value = "  +4.07%  "
$ws.Range("D49").Value = "'4.84"
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("D50").Value = "'94.68"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "'2.500.42"
$ws.Range("E51").Value = "  +0.46%  "
